# Update the "rangoDesde"/"rangoHasta" cheque numbers used by the
# "contraordenar cheque" datadriven test (IOS transaction stabilization).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

$ws.Range("P2").Value = 65785
$ws.Range("Q3").Value = 65772
$ws.Range("R3").Value = 65773
